$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename activity from "Theme Party" to "Kolokium Zon Selatan" (cell C40).
$ws.Range("C40").Value = "Kolokium Zon Selatan"

# 2. Update the merit value for that activity (cell D40): 100 -> 1800.
$ws.Range("D40").Value = 1800

# 3. Merge B15:C15 (label "Penandaan Fail" spans both columns, matching the
#    pattern already used by row 14's B14:C14 merge).
$ws.Range("B15:C15").Merge()

# 4. Give C15 the same (bold, non-wrapping) formatting already used by the
#    sibling sub-total label cells C21/C27/C33, so it matches that look.
$ws.Range("C21").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
